# Atualização automática de SAPIRANGA.xlsx
#
# - Renomeia "Paineis DARQ"            -> "PAINEIS DARQ"
# - Renomeia "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove a planilha "Desarquivamentos Pendentes"

$wb = $excel.ActiveWorkbook

# Rename sheets to their updated (upper-case) titles.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the obsolete "Desarquivamentos Pendentes" sheet.
# Suppress the "permanently delete" confirmation Excel normally shows.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
